$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 451.5116343142244
$ws.Range("C26").Value = 84
$ws.Range("D26").Value = 90374.39999999999
$ws.Range("E26").Value = 306.8912637288979
$ws.Range("F26").Value = 50436
$ws.Range("G26").Value = 54129.6
$ws.Range("H26").Value = 135
$ws.Range("I26").Value = 35
$ws.Range("J26").Value = 2
$ws.Range("K26").Value = 165
$ws.Range("L26").Value = 60
$ws.Range("M26").Value = 28
$ws.Range("N26").Value = 1171
$ws.Range("O26").Value = 0.28
$ws.Range("P26").Value = 2.325
$ws.Range("Q26").Value = 0.015
$ws.Range("R26").Value = 28000
$ws.Range("S26").Value = 315
$ws.Range("T26").Value = 400
$ws.Range("U26").Value = 240
$ws.Range("V26").Value = 90000
$ws.Range("W26").Value = 50000

$ws.Range("J31").Select()
